$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.008.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "'1.556.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'286.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").Value = "'0.3746"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("D8").Value = "'0.3249"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").Value = "'1.128"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").Value = "'41.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -13.42%  "
$ws.Range("D11").Value = "'0.07304"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.49%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "'19.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.24%  "
$ws.Range("D14").Value = "'5.742"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.66%  "
$ws.Range("D15").Value = "'6.852"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "'1.564.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'0.00001082"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.70%  "
$ws.Range("D18").Value = "'0.06635"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "'85.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("D20").Value = "'6.430"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'15.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("D23").Value = "'11.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("D24").Value = "'22.063.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").Value = "'2.259"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.30%  "
$ws.Range("D26").Value = "'2.517"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.61%  "
$ws.Range("D27").Value = "'149.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").Value = "'18.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").Value = "'4.846"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.94%  "
$ws.Range("D30").Value = "'1.740.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "'120.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'1.118"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("D34").Value = "'1.729"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.93%  "
$ws.Range("D35").Value = "'9.270"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.05%  "
$ws.Range("D36").Value = "'0.08103"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.97%  "
$ws.Range("D37").Value = "'5.221"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("D38").Value = "'0.02276"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.61%  "
$ws.Range("D39").Value = "'0.06138"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.11%  "
$ws.Range("D40").Value = "'0.2125"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.18%  "
$ws.Range("D41").Value = "'1.213"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.61%  "
$ws.Range("D42").Value = "'10.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.30%  "
$ws.Range("D43").Value = "'1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "'0.5929"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'13.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.97%  "
$ws.Range("D46").Value = "'3.723"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("E47").Value = "  -6.42%  "
$ws.Range("D48").Value = "'1.952"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.09%  "
$ws.Range("D49").Value = "'119.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.58%  "
$ws.Range("E50").Value = "  -4.76%  "
$ws.Range("D51").Value = "'0.06947"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.79%  "
